$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.160.23'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.10'
$ws.Range('E5').Value = '  +0.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.99'
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +1.82%  '
$ws.Range('E9').Value = '  -1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.75'
$ws.Range('E10').Value = '  -1.05%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.355'
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.147'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.82'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').Value = '2.854.67'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '60.093.31'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '2.442.07'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.25'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.51'
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '327.58'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  +0.98%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.14'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.178'
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.65'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.41'
$ws.Range('E27').Value = '  +4.87%  '
$ws.Range('D28').Value = '0.0₃0773'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.77'
$ws.Range('E29').Value = '  -0.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '170.25'
$ws.Range('E30').Value = '  -0.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.11'
$ws.Range('E31').Value = '  -2.28%  '
$ws.Range('E32').Value = '  +3.73%  '
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.52'
$ws.Range('E34').Value = '  -0.91%  '
$ws.Range('E35').Value = '  +3.26%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.22'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '324.67'
$ws.Range('E39').Value = '  +3.70%  '
$ws.Range('E40').Value = '  -1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '145.95'
$ws.Range('E41').Value = '  +4.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.65'
$ws.Range('E42').Value = '  -1.02%  '
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.76'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.579'
$ws.Range('E46').Value = '  +0.46%  '
$ws.Range('E47').Value = '  -1.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '11.05'
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('E49').Value = '  -1.37%  '
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.941'
$ws.Range('E51').Value = '  -1.31%  '
